# Add three new market test-data sheets (Netherlands, Austria, Denmark),
# cloned from the existing "Greece" sheet template, each placed right
# after the previous sheet in tab order.

$wb = $excel.ActiveWorkbook
$greece = $wb.Worksheets.Item("Greece")

# ---- Netherlands -------------------------------------------------------
$greece.Copy($null, $greece)
$nl = $wb.Worksheets.Item("Greece (2)")
$nl.Name = "Netherlands"

$nl.Range("B4").Value = "NGC-3144/T2176"
$nl.Range("B2").Value = "Netherlands Market"

$nl.Columns.Item(1).ColumnWidth = 23.0833
$nl.Columns.Item(2).ColumnWidth = 15.9167
$nl.Columns.Item(3).ColumnWidth = 12.75
$nl.Columns.Item(4).ColumnWidth = 14.75

$nl.Rows.Item(2).RowHeight = 28.8

$nl.Range("E15").Select()

# ---- Austria -------------------------------------------------------
$greece.Copy($null, $nl)
$at = $wb.Worksheets.Item("Greece (2)")
$at.Name = "Austria"

$at.Range("B4").Value = "NGC-3817/T2272"
$at.Range("B2").Value = "Austria Market"

$at.Columns.Item(1).ColumnWidth = 23.0833
$at.Columns.Item(2).ColumnWidth = 15.9167
$at.Columns.Item(3).ColumnWidth = 12.75
$at.Columns.Item(4).ColumnWidth = 14.75

$at.Range("F22").Select()

# ---- Denmark -------------------------------------------------------
$greece.Copy($null, $at)
$dk = $wb.Worksheets.Item("Greece (2)")
$dk.Name = "Denmark"

$dk.Range("B4").Value = "NGC-2913/T2749"
$dk.Range("B2").Value = "Denmark Market"

$dk.Columns.Item(1).ColumnWidth = 23.0833
$dk.Columns.Item(2).ColumnWidth = 15.9167
$dk.Columns.Item(3).ColumnWidth = 12.75
$dk.Columns.Item(4).ColumnWidth = 14.75

$dk.Range("F22").Select()

# Netherlands ends up the active/selected tab.
$nl.Activate()
$nl.Range("E15").Select()
